# Updated symbol list on Mon Jan 16 19:40:06 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the crypto rows on the active sheet. All cells in these two columns are
# stored as plain text (e.g. "299.68", "-0.42%"), so the script forces a
# text number format ("@") before assigning the new value -- otherwise
# Excel would silently reinterpret a string like "299.68" as a number or
# "-0.42%" as a percentage and change the underlying cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, taken row by row from the diff.
$updates = [ordered]@{
    "D2"  = "299.68";     "E2"  = "-0.42%"
    "D3"  = "31.72";      "E3"  = "1.12%"
    "D4"  = "5.104";      "E4"  = "-0.88%"
    "D5"  = "0.08204";    "E5"  = "11.27%"
    "D6"  = "2.589";      "E6"  = "6.46%"
    "D7"  = "7.836";      "E7"  = "-1.46%"
    "D8"  = "3.840";      "E8"  = "1.37%"
    "D9"  = "0.9275";     "E9"  = "0.75%"
    "D10" = "0.1758";     "E10" = "1.67%"
    "D11" = "0.07485";    "E11" = "-1.87%"
    "D12" = "0.08958";    "E12" = "10.64%"
    "D13" = "0.03005";    "E13" = "-0.73%"
    "D14" = "0.1000";     "E14" = "0.75%"
    "D15" = "0.001510";   "E15" = "0.65%"
    "D16" = "0.005891";   "E16" = "-3.95%"
    "D17" = "3.594";      "E17" = "3.71%"
    "E18" = "1.54%"
    "E19" = "-1.16%"
    "D20" = "0.1346";     "E20" = "0.70%"
    "D21" = "3.900";      "E21" = "-16.16%"
    "D22" = "0.1677";     "E22" = "7.01%"
    "D23" = "0.04605";    "E23" = "-1.03%"
    "D24" = "0.001244";   "E24" = "1.75%"
    "D25" = "0.004546";   "E25" = "1.31%"
    "D26" = "0.0001196";  "E26" = "-7.89%"
    "E27" = "81.63%"
    "D39" = "0.01771";    "E39" = "2.38%"
    "D40" = "0.04542";    "E40" = "0.39%"
    "D41" = "0.006948";   "E41" = "-2.97%"
    "D42" = "0.1379";     "E42" = "2.30%"
    "D43" = "0.002203";   "E43" = "-1.56%"
    "D44" = "0.009588";   "E44" = "-10.52%"
    "D45" = "0.00006195"; "E45" = "-1.26%"
    "E46" = "-0.32%"
    "D48" = "0.8068";     "E48" = "-58.16%"
    "D49" = "0.00002094"; "E49" = "-0.32%"
    "E50" = "-0.25%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
